$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.1022823384
$ws.Range("B2").Value = 0.03515795249999999
$ws.Range("B3").Value = 0.049752832400000016
$ws.Range("B4").Value = 0.062374867224999984
$ws.Range("B5").Value = 0.0675031711
$ws.Range("B6").Value = 0.06984271876666666
$ws.Range("B7").Value = 0.07650126641428573
$ws.Range("B8").Value = 0.08085118126250002
$ws.Range("B9").Value = 0.10109256224444448
